# Update Resume and CV
# Applies four small date-formatting tweaks found in the diff:
#   1. "...     Expected June, 2025 " -> "...       Expected June, 2025 "
#      (two extra leading spaces inserted right before "Expected")
#   2. "...     June, 2020 " -> "...       June 2020 "
#      (comma after "June" dropped, two extra leading spaces added)
#   3. "2020-Present" -> still "2020-Present" visually, but the run is
#      split into "2020" and "-Present" (mirrors the proofing-driven
#      run split in the source edit)
#   4. "...   2020" (ARCS Scholar line) -> "...   <tab>    2020"
#      (an extra tab + four spaces inserted before the trailing "2020")

$d = $word.ActiveDocument

# --- 1) "Expected " gains two leading spaces ------------------------------
$d.Content.Find.Execute("Expected", $true, $false, $false, $false, $false, `
    $true, 1, $false, "  Expected", 2) | Out-Null

# --- 2) "June, 2020" -> "June 2020" with two extra leading spaces --------
$d.Content.Find.Execute("     June, 2020 ", $true, $false, $false, $false, `
    $false, $true, 1, $false, "       June 2020 ", 2) | Out-Null

# --- 3) split "2020-Present" into "2020" + "-Present" ---------------------
$rng = $d.Content
$rng.Find.Execute("2020-Present", $true, $false, $false, $false, $false, `
    $true, 1, $false, "", 0) | Out-Null
$yearStart = $rng.Start
$yearRange = $d.Range($yearStart, $yearStart + 4)
# Nudge the font color off and back to "automatic" on just the "2020"
# portion so Word is forced to materialize it as its own run, the same
# way the original run ended up split in two.
$yearRange.Font.Color = 0
$yearRange.Font.Color = -16777216

# --- 4) extra tab + 4 spaces before the trailing "2020" on ARCS Scholar --
$d.Content.Find.Execute("   2020", $true, $false, $false, $false, $false, `
    $true, 1, $false, ("   " + [char]9 + "    2020"), 2) | Out-Null
